$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections on existing rows (row 2 and row 4 quantity / total) ---
$ws.Range("L2").Value = 1
$ws.Range("N2").Value = 10250

$ws.Range("L4").Value = 5
$ws.Range("N4").Value = 97500

# --- Grow the table (Semaine_1) from A1:P10 to A1:P14 ---
$lo = $ws.ListObjects.Item("Semaine_1")
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

# --- Row 11 ---
$ws.Range("A11").Value = 45924
$ws.Range("B11").Value = "Ndack NDAO"
$ws.Range("C11").Value = "GUEDIAWAYE"
$ws.Range("D11").Value = "Marché Bou Bess"
$ws.Range("E11").Value = "MOUSTAPHA BAKHDAD"
$ws.Range("F11").Value = 776180875
$ws.Range("G11").Value = "Grossiste"
$ws.Range("H11").Value = "Client Partenaire"
$ws.Range("I11").Value = "Livraison"
$ws.Range("J11").Value = "Ok"
$ws.Range("K11").Value = "Café pot Refraish 200g"
$ws.Range("L11").Value = 25
$ws.Range("M11").Value = 19500
$ws.Range("N11").Value = 487500
$ws.Range("O11").Formula = '="S"&_xlfn.ISOWEEKNUM(Semaine_1[[#This Row],[Date]])'
$ws.Range("P11").Formula = '=TEXT(Semaine_1[[#This Row],[Date]],"MMMM")'

# --- Row 12 ---
$ws.Range("A12").Value = 45924
$ws.Range("B12").Value = "Seynabou SOW"
$ws.Range("C12").Value = "CASTOR"
$ws.Range("D12").Value = "Ben Tally"
$ws.Range("E12").Value = "Lye"
$ws.Range("F12").Value = 775426848
$ws.Range("G12").Value = "Demi-Gros"
$ws.Range("H12").Value = "Client Partenaire"
$ws.Range("I12").Value = "Livraison"
$ws.Range("J12").Value = "Merci beaucoup "
$ws.Range("K12").Value = "Café pot Refraish 200g"
$ws.Range("L12").Value = 20
$ws.Range("M12").Value = 19500
$ws.Range("N12").Value = 390000
$ws.Range("O12").Formula = '="S"&_xlfn.ISOWEEKNUM(Semaine_1[[#This Row],[Date]])'
$ws.Range("P12").Formula = '=TEXT(Semaine_1[[#This Row],[Date]],"MMMM")'

# --- Row 13 ---
$ws.Range("A13").Value = 45924
$ws.Range("B13").Value = "Seynabou SOW"
$ws.Range("C13").Value = "CASTOR"
$ws.Range("D13").Value = "Ben Tally"
$ws.Range("E13").Value = "Lye"
$ws.Range("F13").Value = 775426848
$ws.Range("G13").Value = "Demi-Gros"
$ws.Range("H13").Value = "Client Partenaire"
$ws.Range("I13").Value = "Livraison"
$ws.Range("J13").Value = "Merci beaucoup "
$ws.Range("K13").Value = "Café pot Refraish 50g"
$ws.Range("L13").Value = 5
$ws.Range("M13").Value = 10250
$ws.Range("N13").Value = 51250
$ws.Range("O13").Formula = '="S"&_xlfn.ISOWEEKNUM(Semaine_1[[#This Row],[Date]])'
$ws.Range("P13").Formula = '=TEXT(Semaine_1[[#This Row],[Date]],"MMMM")'

# --- Row 14 ---
$ws.Range("A14").Value = 45925
$ws.Range("B14").Value = "DIATTA FAYE"
$ws.Range("C14").Value = "PNR"
$ws.Range("D14").Value = "Bargny"
$ws.Range("E14").Value = "Wakeur Alpha Thiombane"
$ws.Range("F14").Value = 783758073
$ws.Range("G14").Value = "Grossiste"
$ws.Range("H14").Value = "Client Partenaire"
$ws.Range("I14").Value = "Livraison"
$ws.Range("J14").Value = "Il attend son café refraish 1,5 pour demain"
$ws.Range("K14").Value = "Café pot Refraish 200g"
$ws.Range("L14").Value = 25
$ws.Range("M14").Value = 19500
$ws.Range("N14").Value = 487500
$ws.Range("O14").Formula = '="S"&_xlfn.ISOWEEKNUM(Semaine_1[[#This Row],[Date]])'
$ws.Range("P14").Formula = '=TEXT(Semaine_1[[#This Row],[Date]],"MMMM")'

# --- Formatting: new rows use Times New Roman, date format on A, #,##0 on M/N, wrap text on J ---
$ws.Range("A11:A14").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("M11:N14").NumberFormat = "#,##0"

$ws.Range("A11:N14").Font.Name = "TIMES"
$ws.Range("J11:J14").WrapText = $true

# --- Selection matches final saved state ---
$ws.Range("A14:N14").Select()
